$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '21.775.60'
$ws.Range('E2').Value = '  -1.69%  '
$ws.Range('D3').Value = '1.541.92'
$ws.Range('E3').Value = '  -1.23%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '289.45'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3899'
$ws.Range('E7').Value = '  +2.47%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3189'
$ws.Range('E8').Value = '  -2.93%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '43.34'
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07190'
$ws.Range('E10').Value = '  -2.34%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.066'
$ws.Range('E11').Value = '  -6.31%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.627'
$ws.Range('E13').Value = '  -3.51%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '18.61'
$ws.Range('E14').Value = '  -6.80%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.604'
$ws.Range('E15').Value = '  -3.93%  '
$ws.Range('D16').Value = '1.540.57'
$ws.Range('E16').Value = '  -1.52%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001108'
$ws.Range('E17').Value = '  +1.26%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06578'
$ws.Range('E18').Value = '  -0.91%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '83.38'
$ws.Range('E19').Value = '  -2.49%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.9996'
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.146'
$ws.Range('E21').Value = '  -4.95%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '15.38'
$ws.Range('E22').Value = '  -4.66%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.86'
$ws.Range('E23').Value = '  -7.46%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.361'
$ws.Range('E24').Value = '  +4.29%  '
$ws.Range('D25').Value = '21.792.77'
$ws.Range('E25').Value = '  -1.67%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.391'
$ws.Range('E26').Value = '  -5.84%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '145.39'
$ws.Range('E27').Value = '  -3.94%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.42'
$ws.Range('E28').Value = '  -3.53%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.839'
$ws.Range('E29').Value = '  -0.58%  '
$ws.Range('D30').Value = '1.719.24'
$ws.Range('E30').Value = '  -1.15%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '117.76'
$ws.Range('E31').Value = '  -3.04%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.917'
$ws.Range('E32').Value = '  -1.83%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.9684'
$ws.Range('E33').Value = '  -13.73%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.08195'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '8.972'
$ws.Range('E35').Value = '  -4.03%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.06136'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.128'
$ws.Range('E37').Value = '  -3.06%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02211'
$ws.Range('E38').Value = '  -4.14%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.2043'
$ws.Range('E39').Value = '  -4.47%  '
$ws.Range('E40').Value = '  -3.65%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.436'
$ws.Range('E41').Value = '  -23.57%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.9999'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '10.67'
$ws.Range('E43').Value = '  -3.56%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.5763'
$ws.Range('E44').Value = '  -3.67%  '
$ws.Range('E45').Value = '  -4.24%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.737'
$ws.Range('E46').Value = '  -0.64%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5535'
$ws.Range('E47').Value = '  -4.41%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '117.79'
$ws.Range('E48').Value = '  -2.63%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.876'
$ws.Range('E49').Value = '  -5.74%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.132'
$ws.Range('E50').Value = '  -3.37%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06733'
$ws.Range('E51').Value = '  -3.70%  '
